# Scheduled runner update: refresh market-board derived price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job Leve
# profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 255.70589
$ws.Range("J9").Value = 307.5
$ws.Range("L9").Value = 307.5
$ws.Range("N9").Value = -645.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").ClearContents()
$ws.Range("N45").Value = 0

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 19961
$ws.Range("J95").Value = 19961
$ws.Range("L95").Value = 19961
$ws.Range("N95").Value = -25453

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1618.52
$ws.Range("I111").Value = 1023.2857
$ws.Range("K111").Value = 3069.8571
$ws.Range("M111").Value = -2.857100000000173

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 9000
$ws.Range("I116").Value = 3000
$ws.Range("K116").Value = 3000
$ws.Range("M116").Value = 442

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 12434.944
$ws.Range("I132").Value = 14256.267
$ws.Range("J132").Value = 3328.3333
$ws.Range("K132").Value = 42768.801
$ws.Range("L132").Value = 9984.999899999999
$ws.Range("M132").Value = -40238.801
$ws.Range("N132").Value = -15044.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2709.8276
$ws.Range("I137").Value = 1937.6923
$ws.Range("K137").Value = 5813.0769
$ws.Range("M137").Value = -3263.0769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 916.3158
$ws.Range("J2").Value = 995.55554
$ws.Range("L2").Value = 995.55554
$ws.Range("N2").Value = -1221.55554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 46109
$ws.Range("J92").Value = 46109
$ws.Range("L92").Value = 46109
$ws.Range("N92").Value = -51101

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 916.3158
$ws.Range("J116").Value = 995.55554
$ws.Range("L116").Value = 995.55554
$ws.Range("N116").Value = -5583.55554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 34849.5
$ws.Range("J119").Value = 34849.5
$ws.Range("L119").Value = 34849.5
$ws.Range("N119").Value = -44525.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5750
$ws.Range("I132").Value = 6166.6665
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 18499.9995
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -15969.9995
$ws.Range("N132").Value = -18560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 916.3158
$ws.Range("J3").Value = 995.55554
$ws.Range("L3").Value = 995.55554
$ws.Range("N3").Value = -1223.55554

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1166.7693
$ws.Range("I20").Value = 1052.6666
$ws.Range("J20").Value = 1423.5
$ws.Range("K20").Value = 1052.6666
$ws.Range("L20").Value = 1423.5
$ws.Range("M20").Value = -805.6666
$ws.Range("N20").Value = -1917.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3757.7778
$ws.Range("I134").Value = 3757.7778
$ws.Range("K134").Value = 11273.3334
$ws.Range("M134").Value = -8738.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4788.686
$ws.Range("I31").Value = 2801.5
$ws.Range("K31").Value = 2801.5
$ws.Range("M31").Value = -2506.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4788.686
$ws.Range("I34").Value = 2801.5
$ws.Range("K34").Value = 2801.5
$ws.Range("M34").Value = -2599.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2689.8333
$ws.Range("I58").Value = 1036.25
$ws.Range("K58").Value = 1036.25
$ws.Range("M58").Value = -833.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 4020.3333
$ws.Range("I94").Value = 2154.7144
$ws.Range("J94").Value = 6632.2
$ws.Range("K94").Value = 2154.7144
$ws.Range("L94").Value = 6632.2
$ws.Range("M94").Value = -1703.7144
$ws.Range("N94").Value = -7534.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2689.8333
$ws.Range("I136").Value = 1036.25
$ws.Range("K136").Value = 3108.75
$ws.Range("M136").Value = -558.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 793.44446
$ws.Range("I5").Value = 728
$ws.Range("K5").Value = 2184
$ws.Range("M5").Value = -2072

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1303.421
$ws.Range("J68").Value = 1285
$ws.Range("L68").Value = 3855
$ws.Range("N68").Value = -5477

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1303.421
$ws.Range("J71").Value = 1285
$ws.Range("L71").Value = 11565
$ws.Range("N71").Value = -19677

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 380
$ws.Range("J97").Value = 355
$ws.Range("L97").Value = 1065
$ws.Range("N97").Value = -2057

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 4722
$ws.Range("I132").Value = 4444
$ws.Range("K132").Value = 39996
$ws.Range("M132").Value = -37466

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 793.44446
$ws.Range("I135").Value = 728
$ws.Range("K135").Value = 6552
$ws.Range("M135").Value = -4017

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 7565.6665
$ws.Range("J136").Value = 7679
$ws.Range("L136").Value = 23037
$ws.Range("N136").Value = -33237

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2900
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").ClearContents()
$ws.Range("N7").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 2149.8333
$ws.Range("I32").Value = 2149.8333
$ws.Range("K32").Value = 2149.8333
$ws.Range("M32").Value = -1832.8333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("N94").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2540.1428
$ws.Range("I122").Value = 2646.8333
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 7940.499899999999
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -5490.499899999999
$ws.Range("N122").Value = -10600

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("N126").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 16281.286
$ws.Range("J132").Value = 14995.25
$ws.Range("L132").Value = 44985.75
$ws.Range("N132").Value = -50045.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4398.2
$ws.Range("I136").Value = 3264.6667
$ws.Range("J136").Value = 6098.5
$ws.Range("K136").Value = 9794.000100000001
$ws.Range("L136").Value = 18295.5
$ws.Range("M136").Value = -7244.000100000001
$ws.Range("N136").Value = -23395.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 49999
$ws.Range("J119").Value = 49999
$ws.Range("L119").Value = 49999
$ws.Range("N119").Value = -59675

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2908.6667
$ws.Range("I122").Value = 2948.2173
$ws.Range("K122").Value = 8844.651899999999
$ws.Range("M122").Value = -6394.651899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5311.913
$ws.Range("I126").Value = 3101.9
$ws.Range("J126").Value = 7011.923
$ws.Range("K126").Value = 9305.700000000001
$ws.Range("L126").Value = 21035.769
$ws.Range("M126").Value = -6835.700000000001
$ws.Range("N126").Value = -25975.769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1941.6666
$ws.Range("I132").Value = 1941.6666
$ws.Range("K132").Value = 5824.9998
$ws.Range("M132").Value = -3294.9998
